# The workbook contains a single weekly price-log sheet where each row is one
# day's record, sorted by an arbitrary "random" date order. A new weekly
# record was added to the log. In the canonical OOXML this shows up as a new
# row inserted at row 147, with every following row (old 147..235) shifted
# down by one (to 148..236), and the sheet dimension growing from
# A1:R235 to A1:R236.
#
# Reproduce that with a real row insert (so Excel shifts the existing data
# down and copies the row-147 formatting, e.g. the date number format on
# column D) and then populate the newly inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 147; this shifts old rows 147-235 down
# to 148-236 and extends the used range / dimension accordingly.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly record.
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C147").Value = "Los Lagos"
$ws.Range("D147").Value = 44606
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = 100112037
$ws.Range("G147").Value = "Cebollín"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 60
$ws.Range("K147").Value = 5500
$ws.Range("L147").Value = 6000
$ws.Range("M147").Value = 5750
$ws.Range("N147").Value = "$/paquete 36 unidades"
$ws.Range("O147").Value = "Región Metropolitana"
$ws.Range("P147").Value = 160
$ws.Range("Q147").Value = 36
$ws.Range("R147").Value = "Hortaliza"
